$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.559.45'
$ws.Range("E2").Value = '  -2.73%  '
$ws.Range("D3").Value = '3.802.43'
$ws.Range("E3").Value = '  +0.59%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''599.60'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.77%  '
$ws.Range("D6").Value = '''167.90'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.49%  '
$ws.Range("D7").Value = '3.798.57'
$ws.Range("E7").Value = '  +0.51%  '
$ws.Range("D8").Value = '''0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '''0.528'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.14%  '
$ws.Range("D10").Value = '''0.158'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.58%  '
$ws.Range("D11").Value = '''6.18'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.15%  '
$ws.Range("D12").Value = '''0.464'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.56%  '
$ws.Range("D13").Value = '''38.30'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.12%  '
$ws.Range("D14").Value = '''0.0000243'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.75%  '
$ws.Range("D15").Value = '4.437.86'
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").Value = '3.828.49'
$ws.Range("E16").Value = '  +1.30%  '
$ws.Range("D17").Value = '67.624.01'
$ws.Range("E17").Value = '  -2.70%  '
$ws.Range("B18").Value = 'TRON'
$ws.Range("C18").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D18").Value = '''0.115'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.75%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '''7.18'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.06%  '
$ws.Range("D20").Value = '''17.41'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.73%  '
$ws.Range("D21").Value = '''492.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.40%  '
$ws.Range("D22").Value = '''9.36'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.45%  '
$ws.Range("E23").Value = '  +0.62%  '
$ws.Range("D24").Value = '''85.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.04%  '
$ws.Range("E25").Value = '  +2.12%  '
$ws.Range("D26").Value = '''2.37'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.22%  '
$ws.Range("D27").Value = '''12.28'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.80%  '
$ws.Range("D28").Value = '''10.09'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.78%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").Value = '''2.97'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.18%  '
$ws.Range("D31").Value = '''2.42'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.40%  '
$ws.Range("D32").Value = '''32.54'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.73%  '
$ws.Range("D33").Value = '''7.81'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.56%  '
$ws.Range("E34").Value = '  -5.25%  '
$ws.Range("D35").Value = '''1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("E36").Value = '  -3.03%  '
$ws.Range("D37").Value = '''5.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -5.41%  '
$ws.Range("E38").Value = '  -4.82%  '
$ws.Range("D39").Value = '''463.60'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.56%  '
$ws.Range("D40").Value = '''0.327'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.68%  '
$ws.Range("D41").Value = '''49.44'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.84%  '
$ws.Range("D42").Value = '''1.99'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.13%  '
$ws.Range("E43").Value = '  -5.79%  '
$ws.Range("D44").Value = '''8.37'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.55%  '
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").Value = '''40.47'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -8.98%  '
$ws.Range("D47").Value = '2.840.25'
$ws.Range("E47").Value = '  -3.92%  '
$ws.Range("D48").Value = '''140.25'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.00%  '
$ws.Range("D49").Value = '''0.0349'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.83%  '
$ws.Range("D50").Value = '''24.40'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.33%  '
$ws.Range("D51").Value = '''25.71'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.35%  '
